$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1624
$ws.Range("F4").Value = 2116
$ws.Range("F5").Value = 9325
$ws.Range("F7").Value = 1290
$ws.Range("F9").Value = 687
$ws.Range("F10").Value = 612
$ws.Range("F12").Value = 177
$ws.Range("F13").Value = 310
$ws.Range("F15").Value = 63
$ws.Range("F16").Value = 1577
$ws.Range("F17").Value = 1354
$ws.Range("F18").Value = 574
$ws.Range("F20").Value = 1435
$ws.Range("F21").Value = 107
$ws.Range("F22").Value = 272
$ws.Range("F24").Value = 113
$ws.Range("F25").Value = 83
$ws.Range("F26").Value = 77
$ws.Range("F27").Value = 342
$ws.Range("F28").Value = 342
$ws.Range("F29").Value = 1092
$ws.Range("F30").Value = 15
$ws.Range("F31").Value = 41
$ws.Range("F32").Value = 254
$ws.Range("F33").Value = 236
$ws.Range("F34").Value = 68
$ws.Range("F35").Value = 590
$ws.Range("F36").Value = 621
$ws.Range("F40").Value = 170
$ws.Range("F41").Value = 160
$ws.Range("F42").Value = 562
$ws.Range("F43").Value = 1245
$ws.Range("F44").Value = 713
$ws.Range("F45").Value = 257
$ws.Range("F46").Value = 54
$ws.Range("F47").Value = 55

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 164
$ws.Range("F8").Value = 10
$ws.Range("F9").Value = 5
$ws.Range("F11").Value = 680
$ws.Range("F14").Value = 12
$ws.Range("F16").Value = 14
$ws.Range("F19").Value = 954
$ws.Range("F22").Value = 259
$ws.Range("F23").Value = 648
$ws.Range("F25").Value = 275
$ws.Range("F26").Value = 275
$ws.Range("F31").Value = 118
$ws.Range("F33").Value = 28
$ws.Range("F36").Value = 105

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 347
$ws.Range("F7").Value = 2296
$ws.Range("F8").Value = 3423

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1624
$ws.Range("F6").Value = 9325
$ws.Range("F7").Value = 347
$ws.Range("F8").Value = 3423
$ws.Range("F9").Value = 687
$ws.Range("F10").Value = 10
$ws.Range("F11").Value = 5
$ws.Range("F13").Value = 1577
$ws.Range("F14").Value = 680
$ws.Range("F15").Value = 1354
$ws.Range("F18").Value = 1435
$ws.Range("F19").Value = 107
$ws.Range("F20").Value = 272
$ws.Range("F21").Value = 113
$ws.Range("F22").Value = 77
$ws.Range("F23").Value = 342
$ws.Range("F24").Value = 15
$ws.Range("F25").Value = 14
$ws.Range("F27").Value = 41
$ws.Range("F29").Value = 955
$ws.Range("F32").Value = 259
$ws.Range("F33").Value = 590
$ws.Range("F34").Value = 621
$ws.Range("F37").Value = 275
$ws.Range("F39").Value = 562
$ws.Range("F40").Value = 713
$ws.Range("F43").Value = 118
$ws.Range("F45").Value = 28
$ws.Range("F46").Value = 105
$ws.Range("F47").Value = 54
$ws.Range("F48").Value = 55

